$d = $word.ActiveDocument
$t = $d.Tables.Item(2)
$cell = $t.Cell(4, 2)
$cell.Range.Text = "Пламена Георгиева"

$t2 = $d.Tables.Item(2)
$cell2 = $t2.Cell(4, 2)
$ins = $cell2.Range
$ins.Font.Name = "Calibri"
$ins.Font.NameFarEast = "Calibri"
$ins.Font.NameBi = "Calibri"
Write-Output "set fonts ok"
$ins.LanguageID = 1026
Write-Output "set lang ok"
